$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author removed two rows from the dataset (the ones labelled "c" and
# "pc" in column A) and appended two new rows ("zy_l" and "zy_r") at the
# bottom, so the sheet stays 48 data rows (A1:D48) but with different
# contents shifted into place.

# Row 11 (label "c") is deleted first; every row below it shifts up by one.
$ws.Rows.Item(11).Delete()

# After that shift, the row that used to be row 24 (label "pc") is now at
# row 23; delete it too, shifting everything below up by one more.
$ws.Rows.Item(23).Delete()

# Append the two new rows at the bottom of the table (rows 47-48).
$ws.Range("A47").Value = "zy_l"
$ws.Range("B47").Value = 61.55
$ws.Range("C47").Value = 5.82
$ws.Range("D47").Value = 61.16

$ws.Range("A48").Value = "zy_r"
$ws.Range("B48").Value = -56.68
$ws.Range("C48").Value = 6.54
$ws.Range("D48").Value = 70.86

# Leave the selection on the last cell that was filled in, matching the
# cursor position the author would have ended on.
[void]$ws.Range("D48").Select()
